# Revert "Updates to new model policy schedule"
# This reverts the workbook content back to the pre-update text/layout.

$wb  = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsLever = $wb.Worksheets.Item("BAEPAbCiPC")

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------

# A1 header text changed from "... Prices ..." (same shared string, unchanged visible text
# actually differs only by row range below) - the title stays the same text.
$wsAbout.Range("A1").Value = "BAEPAbCiPC Boolean Are Energy Prices Affected by Changes in Production Costs"

# A8/A9 unchanged body text (kept verbatim, just noting for completeness)
$wsAbout.Range("A8").Value = "This control lever sets whether energy prices are affected by production costs"
$wsAbout.Range("A9").Value = "in the modeled region."

# Replace the explanatory paragraph (EPS 3.1.0 -> EPS 1.5.0 wording) and shrink it from
# 3 lines + a separate 4-line paragraph down to a single unified 4-line paragraph.
$wsAbout.Range("A11").Value = "As of EPS 1.5.0, this lever supports the three energy carriers (electricity,"
$wsAbout.Range("A12").Value = "district heat, and hydrogen), which tend to be produced and consumed locally."
$wsAbout.Range("A13").Value = "It does not affect other fuel types, whose prices are often determined or influenced"
$wsAbout.Range("A14").Value = "by global markets, so domestic producers' costs are less relevant."

# Drop the old trailing paragraph (rows 15-18) that used to follow.
$wsAbout.Rows("15:18").Delete()

# ---------------------------------------------------------------------------
# Sheet "BAEPAbCiPC"
# ---------------------------------------------------------------------------

# Row 1: label simplified from "Unit: boolean (1 or 0)" to "Boolean", and drop the
# italic styling it used to have (now plain).
$wsLever.Range("A1").Value = "Boolean"
$wsLever.Range("A1").ClearFormats()

# Carrier rows (electricity / heat / hydrogen) switch from the old "applyFill,no-fill"
# style to a green highlight fill.
$wsLever.Range("A2:B2").Interior.Color = 5296274
$wsLever.Range("A15:B15").Interior.Color = 5296274
$wsLever.Range("A22:B22").Interior.Color = 5296274

# "nuclear" and "municipal solid waste" rows lose their "(NOT USED)" suffix and their
# grey highlighting - they become ordinary, unstyled fuel rows.
$wsLever.Range("A5").Value = "nuclear"
$wsLever.Range("A5:B5").ClearFormats()
$wsLever.Range("A21").Value = "municipal solid waste"
$wsLever.Range("A21:B21").ClearFormats()

Write-Host "done"
